$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 803
$ws.Range("F3").Value = 2876
$ws.Range("F4").Value = 1349
$ws.Range("F5").Value = 1963
$ws.Range("F10").Value = 89
$ws.Range("F11").Value = 11955
$ws.Range("F12").Value = 6760
$ws.Range("F15").Value = 430
$ws.Range("F19").Value = 942
$ws.Range("F20").Value = 102
$ws.Range("F22").Value = 944
$ws.Range("F23").Value = 3674
$ws.Range("F24").Value = 61
$ws.Range("F25").Value = 993
$ws.Range("F32").Value = 42
$ws.Range("F33").Value = 317
$ws.Range("F34").Value = 5063
$ws.Range("F36").Value = 1267
$ws.Range("F37").Value = 250
$ws.Range("F38").Value = 709
$ws.Range("F39").Value = 1227
$ws.Range("F40").Value = 553

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 3704
$ws.Range("F15").Value = 18
$ws.Range("F16").Value = 3
$ws.Range("F25").Value = 40

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 9120
$ws.Range("F3").Value = 518
$ws.Range("F4").Value = 1865

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 9120
$ws.Range("F3").Value = 518
$ws.Range("F4").Value = 1865
$ws.Range("F5").Value = 803
$ws.Range("F6").Value = 2876
$ws.Range("F9").Value = 1349
$ws.Range("F13").Value = 89
$ws.Range("F14").Value = 11955
$ws.Range("F15").Value = 6760
$ws.Range("F17").Value = 3704
$ws.Range("F20").Value = 430
$ws.Range("F23").Value = 942
$ws.Range("F24").Value = 102
$ws.Range("F26").Value = 944
$ws.Range("F27").Value = 3674
$ws.Range("F28").Value = 61
$ws.Range("F29").Value = 993
$ws.Range("F37").Value = 3
$ws.Range("F38").Value = 317
$ws.Range("F39").Value = 1267
$ws.Range("F40").Value = 250
$ws.Range("F42").Value = 1227
$ws.Range("F43").Value = 553
